$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.915.11'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '1.891.96'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.45'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4585'
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3901'
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07850'
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9892'
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.88'
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("D12").Value = '1.943.70'
$ws.Range("E12").Value = '  +9.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.042'
$ws.Range("E13").Value = '  +2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.691'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06946'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.06'
$ws.Range("E16").Value = '  +1.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009982'
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").Value = '28.922.09'
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.297'
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("D24").Value = '2.118.91'
$ws.Range("E24").Value = '  +5.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.059'
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.06'
$ws.Range("E26").Value = '  +1.79%  '
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.923'
$ws.Range("E28").Value = '  +4.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.930'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.49'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09358'
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9083'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.294'
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.332'
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.260'
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.191'
$ws.Range("E36").Value = '  +2.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05770'
$ws.Range("E37").Value = '  +1.47%  '
$ws.Range("E38").Value = '  +1.39%  '
$ws.Range("E39").Value = '  -0.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.733'
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5682'
$ws.Range("E41").Value = '  +2.42%  '
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.748'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.300'
$ws.Range("E44").Value = '  +9.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.94'
$ws.Range("E45").Value = '  +3.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5355'
$ws.Range("E46").Value = '  +2.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07043'
$ws.Range("E47").Value = '  -1.07%  '
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.83'
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.526'
$ws.Range("E50").Value = '  +3.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.065'
$ws.Range("E51").Value = '  -5.33%  '
